# Update countries & provincias Spain
# Applies the 24-Jun-2020 04:49 data refresh to the "Pais" sheet:
#   - refreshed case counts for Bolivia, Haiti and Nueva Zelanda
#   - the table is kept sorted by "Casos totales" (col B) descending, so a
#     few countries swap places with their neighbour once the numbers move
#     (Haiti/Macedonia, Dominica/Fiyi, Groenlandia/Islas Malvinas,
#     Seychelles/Montserrat)
#   - the "datos actualizados" timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($r, 1).Value = $country
    $ws.Cells.Item($r, 2).Value = $total
    $ws.Cells.Item($r, 3).Value = $nuevos
    $ws.Cells.Item($r, 4).Value = $activos
    $ws.Cells.Item($r, 5).Value = $recuperados
    $ws.Cells.Item($r, 6).Value = $criticos
    $ws.Cells.Item($r, 7).Value = $muertesHoy
    $ws.Cells.Item($r, 8).Value = $muertes
}

# Bolivia: updated numbers (row stays put, no-one overtakes it)
Set-Row 47 "Bolivia" 26389 896 6300 19243 0 26 846

# Haiti jumps ahead of Republica de Macedonia after its update
Set-Row 80 "Haiti" 5324 113 436 4799 0 1 89
Set-Row 81 "Republica de Macedonia" 5311 0 2048 3012 0 0 251

# Nueva Zelanda: updated numbers (row stays put)
Set-Row 119 "Nueva Zelanda" 1516 1 1483 11 0 0 22

# Dominica / Fiyi swap order (tied totals)
Set-Row 202 "Dominica" 18 0 18 0 0 0 0
Set-Row 203 "Fiyi" 18 0 18 0 0 0 0

# Groenlandia / Islas Malvinas swap order (tied totals)
Set-Row 208 "Groenlandia" 13 0 13 0 0 0 0
Set-Row 209 "Islas Malvinas" 13 0 13 0 0 0 0

# Seychelles / Montserrat swap order (tied totals)
Set-Row 211 "Seychelles" 11 0 11 0 0 0 0
Set-Row 212 "Montserrat" 11 0 10 0 0 0 1

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 04:49"
